$wb = $excel.ActiveWorkbook

# Update the zh-cn sheet: row 2 handoff/handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 04:36:05"
$wsZhCn.Range("H2").Value = "2016-03-20 04:37:02"

# Update the de-de sheet: row 2 handoff/handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 04:36:13"
$wsDeDe.Range("H2").Value = "2016-03-20 04:37:16"
